# Insert a new data row at row 683 (pushing the existing row 683..740 down to
# 684..741), then populate the new row with the latest weekly price record.
# This mirrors the upstream diff: dimension grows from A1:R740 to A1:R741 and
# every row from 683 onward is the previous week's row shifted down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("683:683").Insert()

$ws.Range("A683").Value = 6
$ws.Range("B683").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C683").Value = 'Metropolitana'
$ws.Range("D683").Value = 45223
$ws.Range("E683").Value = 13
$ws.Range("F683").Value = 100112030
$ws.Range("G683").Value = 'Poroto granado'
$ws.Range("H683").Value = 'Sin especificar'
$ws.Range("I683").Value = 'Primera'
$ws.Range("J683").Value = 230
$ws.Range("K683").Value = 45000
$ws.Range("L683").Value = 45000
$ws.Range("M683").Value = 45000
$ws.Range("N683").Value = '$/malla 25 kilos'
$ws.Range("O683").Value = 'Perú'
$ws.Range("P683").Value = 1800
$ws.Range("Q683").Value = 25
$ws.Range("R683").Value = 'Hortaliza'
